# Update "cryptos" worksheet with refreshed price data (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) is stored as text (numeric-looking strings such as
# "0.06200" that must keep their exact formatting/trailing zeros). Force the
# whole column range to Text format before writing so Excel does not
# auto-convert the assigned strings into real numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Price updates -------------------------------------------------------
$ws.Range("D2").Value = "267.33"
$ws.Range("D3").Value = "21.62"
$ws.Range("D4").Value = "6.283"
$ws.Range("D5").Value = "0.06200"
$ws.Range("D6").Value = "3.571"
$ws.Range("D7").Value = "6.535"
$ws.Range("D8").Value = "1.394"
$ws.Range("D9").Value = "0.8239"
$ws.Range("D10").Value = "0.1629"
$ws.Range("D11").Value = "0.08224"
$ws.Range("D12").Value = "0.03568"
$ws.Range("D13").Value = "0.03208"
$ws.Range("D14").Value = "0.09197"
$ws.Range("D15").Value = "3.773"
$ws.Range("D16").Value = "0.001642"
$ws.Range("D17").Value = "0.04656"
$ws.Range("D18").Value = "0.006361"
$ws.Range("D19").Value = "0.006184"
$ws.Range("D21").Value = "0.0001501"
$ws.Range("D23").Value = "2.237"
$ws.Range("D25").Value = "0.3338"
$ws.Range("D28").Value = "0.0002714"
$ws.Range("D40").Value = "0.04702"
$ws.Range("D41").Value = "0.006927"

# --- Rows 42 & 43 swapped places in the ranking (symbol list reorder) ----
# Row 42 becomes BKEXToken, row 43 becomes CEJI, with new data for each.
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1119"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.002902"
$ws.Range("E43").Value = "42CEJICEJI"

# --- Remaining price updates ---------------------------------------------
$ws.Range("D44").Value = "0.01099"
$ws.Range("D45").Value = "0.00006112"
$ws.Range("D46").Value = "0.0009905"
$ws.Range("D48").Value = "0.9805"
$ws.Range("D49").Value = "0.001137"
$ws.Range("D50").Value = "0.00001901"
$ws.Range("D51").Value = "0.01241"
